$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new date row labels in column A (rows 139-144)
$ws.Cells.Item(139, 1).Value = "17 06 2020"
$ws.Cells.Item(140, 1).Value = "18 06 2020"
$ws.Cells.Item(141, 1).Value = "19 06 2020"
$ws.Cells.Item(142, 1).Value = "20 06 2020"
$ws.Cells.Item(143, 1).Value = "21 06 2020"
$ws.Cells.Item(144, 1).Value = "22 06 2020"

# Fill per-state numeric data for rows 134-139
# Row 134
$ws.Cells.Item(134, 2).Value = 51.70431659
$ws.Cells.Item(134, 3).Value = 4.69340585
$ws.Cells.Item(134, 4).Value = 19.81901565
$ws.Cells.Item(134, 6).Value = 21.41013262
$ws.Cells.Item(134, 7).Value = 5.49657888
$ws.Cells.Item(134, 8).Value = 15.55046069
$ws.Cells.Item(134, 9).Value = 19.87018943
$ws.Cells.Item(134, 10).Value = 15.04428631
$ws.Cells.Item(134, 11).Value = 0
$ws.Cells.Item(134, 12).Value = 9.046151630000001
$ws.Cells.Item(134, 13).Value = 8.1088351
$ws.Cells.Item(134, 15).Value = 6.82665842
$ws.Cells.Item(134, 16).Value = 2.79281869
$ws.Cells.Item(134, 17).Value = 0.24394001
$ws.Cells.Item(134, 18).Value = 7.37937308
$ws.Cells.Item(134, 19).Value = 7.22508909
$ws.Cells.Item(134, 20).Value = 0
$ws.Cells.Item(134, 21).Value = 10.26990046
$ws.Cells.Item(134, 22).Value = 17.42657623
$ws.Cells.Item(134, 23).Value = 16.67125231
$ws.Cells.Item(134, 24).Value = 7.20383191
$ws.Cells.Item(134, 25).Value = 0
$ws.Cells.Item(134, 26).Value = 8.89232056
$ws.Cells.Item(134, 27).Value = 10.4179613
$ws.Cells.Item(134, 28).Value = 7.66490779
$ws.Cells.Item(134, 30).Value = 21.70689912
$ws.Cells.Item(134, 31).Value = 0
$ws.Cells.Item(134, 32).Value = 8.31088928
$ws.Cells.Item(134, 33).Value = 0
$ws.Cells.Item(134, 34).Value = 0
$ws.Cells.Item(134, 35).Value = 14.48860021
$ws.Cells.Item(134, 36).Value = 10.49924398
$ws.Cells.Item(134, 37).Value = 3.42908838
$ws.Cells.Item(134, 38).Value = 6.47732484
$ws.Cells.Item(134, 39).Value = 6.10131629
$ws.Cells.Item(134, 40).Value = 4.00821148
$ws.Cells.Item(134, 41).Value = 15.03714918
$ws.Cells.Item(134, 42).Value = 7.89556947
$ws.Cells.Item(134, 43).Value = 11.21904771
$ws.Cells.Item(134, 45).Value = 3.7793864
$ws.Cells.Item(134, 46).Value = 14.7661431
$ws.Cells.Item(134, 47).Value = 38.02402435
$ws.Cells.Item(134, 48).Value = 9.22121288
$ws.Cells.Item(134, 49).Value = 10.00786153
$ws.Cells.Item(134, 50).Value = 12.18631262
$ws.Cells.Item(134, 51).Value = 9.30897
$ws.Cells.Item(134, 53).Value = 0
$ws.Cells.Item(134, 54).Value = 1.58268273
$ws.Cells.Item(134, 55).Value = 6.70113936
$ws.Cells.Item(134, 56).Value = 12.90796437
$ws.Cells.Item(134, 57).Value = 0

# Row 135
$ws.Cells.Item(135, 2).Value = 42.13616601
$ws.Cells.Item(135, 3).Value = 2.08872188
$ws.Cells.Item(135, 4).Value = 14.77759727
$ws.Cells.Item(135, 6).Value = 20.61961183
$ws.Cells.Item(135, 7).Value = 6.72975819
$ws.Cells.Item(135, 8).Value = 14.6761495
$ws.Cells.Item(135, 9).Value = 21.69454794
$ws.Cells.Item(135, 10).Value = 64.08657862
$ws.Cells.Item(135, 11).Value = 0
$ws.Cells.Item(135, 12).Value = 11.32387081
$ws.Cells.Item(135, 13).Value = 8.3727047
$ws.Cells.Item(135, 15).Value = 4.50517852
$ws.Cells.Item(135, 16).Value = 0.70240984
$ws.Cells.Item(135, 17).Value = 0
$ws.Cells.Item(135, 18).Value = 11.02245303
$ws.Cells.Item(135, 19).Value = 5.83595685
$ws.Cells.Item(135, 20).Value = 0
$ws.Cells.Item(135, 21).Value = 7.19118985
$ws.Cells.Item(135, 22).Value = 21.40604538
$ws.Cells.Item(135, 23).Value = 15.2197034
$ws.Cells.Item(135, 24).Value = 15.86690706
$ws.Cells.Item(135, 25).Value = 0
$ws.Cells.Item(135, 26).Value = 11.20449185
$ws.Cells.Item(135, 27).Value = 15.77574508
$ws.Cells.Item(135, 28).Value = 15.96441776
$ws.Cells.Item(135, 30).Value = 16.87711882
$ws.Cells.Item(135, 31).Value = 0
$ws.Cells.Item(135, 32).Value = 5.83736239
$ws.Cells.Item(135, 33).Value = 0
$ws.Cells.Item(135, 34).Value = 13.25756579
$ws.Cells.Item(135, 35).Value = 10.06045671
$ws.Cells.Item(135, 36).Value = 12.06175255
$ws.Cells.Item(135, 37).Value = 1.82312556
$ws.Cells.Item(135, 38).Value = 4.34633875
$ws.Cells.Item(135, 39).Value = 7.15165137
$ws.Cells.Item(135, 40).Value = 2.4248385
$ws.Cells.Item(135, 41).Value = 27.18952414
$ws.Cells.Item(135, 42).Value = 6.2006259
$ws.Cells.Item(135, 43).Value = 10.78794228
$ws.Cells.Item(135, 45).Value = 2.14056167
$ws.Cells.Item(135, 46).Value = 10.8796859
$ws.Cells.Item(135, 47).Value = 33.06988131
$ws.Cells.Item(135, 48).Value = 9.57764017
$ws.Cells.Item(135, 49).Value = 12.76660331
$ws.Cells.Item(135, 50).Value = 9.31834082
$ws.Cells.Item(135, 51).Value = 6.5834046
$ws.Cells.Item(135, 53).Value = 0
$ws.Cells.Item(135, 54).Value = 2.03125337
$ws.Cells.Item(135, 55).Value = 4.73642319
$ws.Cells.Item(135, 56).Value = 7.79613421
$ws.Cells.Item(135, 57).Value = 0

# Row 136
$ws.Cells.Item(136, 2).Value = 32.65836362
$ws.Cells.Item(136, 3).Value = 3.6103201
$ws.Cells.Item(136, 4).Value = 10.25757835
$ws.Cells.Item(136, 6).Value = 25.77337444
$ws.Cells.Item(136, 7).Value = 6.98822937
$ws.Cells.Item(136, 8).Value = 11.33121216
$ws.Cells.Item(136, 9).Value = 16.06613735
$ws.Cells.Item(136, 10).Value = 53.05671889
$ws.Cells.Item(136, 11).Value = 0
$ws.Cells.Item(136, 12).Value = 13.28040872
$ws.Cells.Item(136, 13).Value = 11.02390005
$ws.Cells.Item(136, 15).Value = 2.57392181
$ws.Cells.Item(136, 16).Value = 0
$ws.Cells.Item(136, 17).Value = 0
$ws.Cells.Item(136, 18).Value = 12.70284185
$ws.Cells.Item(136, 19).Value = 9.619901390000001
$ws.Cells.Item(136, 20).Value = 0
$ws.Cells.Item(136, 21).Value = 4.56329295
$ws.Cells.Item(136, 22).Value = 16.79212632
$ws.Cells.Item(136, 23).Value = 11.82540187
$ws.Cells.Item(136, 24).Value = 19.21001413
$ws.Cells.Item(136, 25).Value = 0
$ws.Cells.Item(136, 26).Value = 11.72202677
$ws.Cells.Item(136, 27).Value = 22.99573928
$ws.Cells.Item(136, 28).Value = 12.82712867
$ws.Cells.Item(136, 30).Value = 12.38611894
$ws.Cells.Item(136, 31).Value = 15.30555642
$ws.Cells.Item(136, 32).Value = 7.66026746
$ws.Cells.Item(136, 33).Value = 0
$ws.Cells.Item(136, 34).Value = 11.04476053
$ws.Cells.Item(136, 35).Value = 6.20057569
$ws.Cells.Item(136, 36).Value = 9.162452200000001
$ws.Cells.Item(136, 37).Value = 15.66837231
$ws.Cells.Item(136, 38).Value = 7.00803294
$ws.Cells.Item(136, 39).Value = 8.788713019999999
$ws.Cells.Item(136, 40).Value = 3.57988818
$ws.Cells.Item(136, 41).Value = 28.05443975
$ws.Cells.Item(136, 42).Value = 4.59849623
$ws.Cells.Item(136, 43).Value = 10.33431854
$ws.Cells.Item(136, 45).Value = 0.80483169
$ws.Cells.Item(136, 46).Value = 7.40734924
$ws.Cells.Item(136, 47).Value = 27.78786411
$ws.Cells.Item(136, 48).Value = 7.53140022
$ws.Cells.Item(136, 49).Value = 13.9720067
$ws.Cells.Item(136, 50).Value = 6.7842094
$ws.Cells.Item(136, 51).Value = 7.31203225
$ws.Cells.Item(136, 53).Value = 0
$ws.Cells.Item(136, 54).Value = 4.39228728
$ws.Cells.Item(136, 55).Value = 11.14946749
$ws.Cells.Item(136, 56).Value = 3.61934772
$ws.Cells.Item(136, 57).Value = 0

# Row 137
$ws.Cells.Item(137, 2).Value = 23.74409981
$ws.Cells.Item(137, 3).Value = 1.7554282
$ws.Cells.Item(137, 4).Value = 17.34767707
$ws.Cells.Item(137, 6).Value = 27.4016708
$ws.Cells.Item(137, 7).Value = 6.9421677
$ws.Cells.Item(137, 8).Value = 8.243027039999999
$ws.Cells.Item(137, 9).Value = 11.0839135
$ws.Cells.Item(137, 10).Value = 42.32705769
$ws.Cells.Item(137, 11).Value = 0
$ws.Cells.Item(137, 12).Value = 13.76575957
$ws.Cells.Item(137, 13).Value = 11.70565977
$ws.Cells.Item(137, 15).Value = 1.02884108
$ws.Cells.Item(137, 16).Value = 7.30955046
$ws.Cells.Item(137, 17).Value = 15.4088639
$ws.Cells.Item(137, 18).Value = 10.7229255
$ws.Cells.Item(137, 19).Value = 7.73414474
$ws.Cells.Item(137, 20).Value = 0
$ws.Cells.Item(137, 21).Value = 9.19166817
$ws.Cells.Item(137, 22).Value = 12.50428827
$ws.Cells.Item(137, 23).Value = 14.05999118
$ws.Cells.Item(137, 24).Value = 15.6527237
$ws.Cells.Item(137, 25).Value = 0
$ws.Cells.Item(137, 26).Value = 12.94118804
$ws.Cells.Item(137, 27).Value = 20.71240431
$ws.Cells.Item(137, 28).Value = 9.853996240000001
$ws.Cells.Item(137, 30).Value = 8.38425559
$ws.Cells.Item(137, 31).Value = 13.59770492
$ws.Cells.Item(137, 32).Value = 8.79379846
$ws.Cells.Item(137, 33).Value = 10.13587195
$ws.Cells.Item(137, 34).Value = 8.90023736
$ws.Cells.Item(137, 35).Value = 21.60315092
$ws.Cells.Item(137, 36).Value = 8.82758318
$ws.Cells.Item(137, 37).Value = 12.8554084
$ws.Cells.Item(137, 38).Value = 12.88196201
$ws.Cells.Item(137, 39).Value = 9.01693796
$ws.Cells.Item(137, 40).Value = 3.33646066
$ws.Cells.Item(137, 41).Value = 31.48444179
$ws.Cells.Item(137, 42).Value = 5.9702498
$ws.Cells.Item(137, 43).Value = 9.613640029999999
$ws.Cells.Item(137, 45).Value = 21.01414392
$ws.Cells.Item(137, 46).Value = 7.12622169
$ws.Cells.Item(137, 47).Value = 22.40339103
$ws.Cells.Item(137, 48).Value = 11.56804438
$ws.Cells.Item(137, 49).Value = 15.06755626
$ws.Cells.Item(137, 50).Value = 8.46708718
$ws.Cells.Item(137, 51).Value = 6.24893827
$ws.Cells.Item(137, 53).Value = 0
$ws.Cells.Item(137, 54).Value = 6.15354995
$ws.Cells.Item(137, 55).Value = 20.32779003
$ws.Cells.Item(137, 56).Value = 0.38493623
$ws.Cells.Item(137, 57).Value = 0

# Row 138
$ws.Cells.Item(138, 2).Value = 15.77041573
$ws.Cells.Item(138, 3).Value = 6.38167513
$ws.Cells.Item(138, 4).Value = 23.71020169
$ws.Cells.Item(138, 6).Value = 26.51272462
$ws.Cells.Item(138, 7).Value = 7.40757907
$ws.Cells.Item(138, 8).Value = 9.420704300000001
$ws.Cells.Item(138, 9).Value = 6.83980607
$ws.Cells.Item(138, 10).Value = 32.21023785
$ws.Cells.Item(138, 11).Value = 0
$ws.Cells.Item(138, 12).Value = 14.28591495
$ws.Cells.Item(138, 13).Value = 11.08580303
$ws.Cells.Item(138, 15).Value = 0
$ws.Cells.Item(138, 16).Value = 14.67345876
$ws.Cells.Item(138, 17).Value = 13.26997064
$ws.Cells.Item(138, 18).Value = 9.7365651
$ws.Cells.Item(138, 19).Value = 5.91243506
$ws.Cells.Item(138, 20).Value = 0
$ws.Cells.Item(138, 21).Value = 13.4346547
$ws.Cells.Item(138, 22).Value = 11.99935221
$ws.Cells.Item(138, 23).Value = 16.6380782
$ws.Cells.Item(138, 24).Value = 15.72519331
$ws.Cells.Item(138, 25).Value = 0
$ws.Cells.Item(138, 26).Value = 9.99906316
$ws.Cells.Item(138, 27).Value = 22.68294042
$ws.Cells.Item(138, 28).Value = 7.12263309
$ws.Cells.Item(138, 30).Value = 17.07461943
$ws.Cells.Item(138, 31).Value = 11.62450955
$ws.Cells.Item(138, 32).Value = 10.25059847
$ws.Cells.Item(138, 33).Value = 7.61638963
$ws.Cells.Item(138, 34).Value = 6.85922549
$ws.Cells.Item(138, 35).Value = 16.8069376
$ws.Cells.Item(138, 36).Value = 10.74776547
$ws.Cells.Item(138, 37).Value = 30.40885957
$ws.Cells.Item(138, 38).Value = 10.1673603
$ws.Cells.Item(138, 39).Value = 9.44845267
$ws.Cells.Item(138, 40).Value = 2.12614932
$ws.Cells.Item(138, 41).Value = 41.86124212
$ws.Cells.Item(138, 42).Value = 9.92145983
$ws.Cells.Item(138, 43).Value = 10.17048866
$ws.Cells.Item(138, 45).Value = 17.67455181
$ws.Cells.Item(138, 46).Value = 16.21495703
$ws.Cells.Item(138, 47).Value = 35.37237619
$ws.Cells.Item(138, 48).Value = 9.097915929999999
$ws.Cells.Item(138, 49).Value = 17.45870941
$ws.Cells.Item(138, 50).Value = 17.94682038
$ws.Cells.Item(138, 51).Value = 6.60244982
$ws.Cells.Item(138, 53).Value = 0
$ws.Cells.Item(138, 54).Value = 4.70917041
$ws.Cells.Item(138, 55).Value = 26.92472059
$ws.Cells.Item(138, 56).Value = 0
$ws.Cells.Item(138, 57).Value = 0

# Row 139
$ws.Cells.Item(139, 2).Value = 8.995002939999999
$ws.Cells.Item(139, 3).Value = 19.95484601
$ws.Cells.Item(139, 4).Value = 29.78224315
$ws.Cells.Item(139, 6).Value = 30.9799694
$ws.Cells.Item(139, 7).Value = 7.15260556
$ws.Cells.Item(139, 8).Value = 6.6276003
$ws.Cells.Item(139, 9).Value = 15.74717158
$ws.Cells.Item(139, 10).Value = 22.98935872
$ws.Cells.Item(139, 11).Value = 26.21685683
$ws.Cells.Item(139, 12).Value = 13.20190782
$ws.Cells.Item(139, 13).Value = 13.77585244
$ws.Cells.Item(139, 15).Value = 0
$ws.Cells.Item(139, 16).Value = 23.96417991
$ws.Cells.Item(139, 17).Value = 11.00146451
$ws.Cells.Item(139, 18).Value = 7.96484606
$ws.Cells.Item(139, 19).Value = 10.99935895
$ws.Cells.Item(139, 20).Value = 0
$ws.Cells.Item(139, 21).Value = 10.43986587
$ws.Cells.Item(139, 22).Value = 8.26607347
$ws.Cells.Item(139, 23).Value = 14.32804681
$ws.Cells.Item(139, 24).Value = 13.84989237
$ws.Cells.Item(139, 25).Value = 0
$ws.Cells.Item(139, 26).Value = 8.61696733
$ws.Cells.Item(139, 27).Value = 19.90160216
$ws.Cells.Item(139, 28).Value = 7.00341663
$ws.Cells.Item(139, 30).Value = 37.71430179
$ws.Cells.Item(139, 31).Value = 37.88782251
$ws.Cells.Item(139, 32).Value = 11.19011656
$ws.Cells.Item(139, 33).Value = 5.50594411
$ws.Cells.Item(139, 34).Value = 18.97673376
$ws.Cells.Item(139, 35).Value = 12.4357182
$ws.Cells.Item(139, 36).Value = 11.46173419
$ws.Cells.Item(139, 37).Value = 25.4200136
$ws.Cells.Item(139, 38).Value = 7.64000837
$ws.Cells.Item(139, 39).Value = 9.510178310000001
$ws.Cells.Item(139, 40).Value = 3.33752991
$ws.Cells.Item(139, 41).Value = 38.37336103
$ws.Cells.Item(139, 42).Value = 13.67848488
$ws.Cells.Item(139, 43).Value = 9.42713382
$ws.Cells.Item(139, 45).Value = 14.31806794
$ws.Cells.Item(139, 46).Value = 15.31797972
$ws.Cells.Item(139, 47).Value = 28.26860419
$ws.Cells.Item(139, 48).Value = 10.73419307
$ws.Cells.Item(139, 49).Value = 17.15001446
$ws.Cells.Item(139, 50).Value = 18.49054225
$ws.Cells.Item(139, 51).Value = 6.99744458
$ws.Cells.Item(139, 53).Value = 17.27835865
$ws.Cells.Item(139, 54).Value = 6.57738704
$ws.Cells.Item(139, 55).Value = 29.19185408
$ws.Cells.Item(139, 56).Value = 0
$ws.Cells.Item(139, 57).Value = 0
